$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 4 (shifts Gearing/CLICK row and below down by one)
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 2).Value = "WAIT"

# Update old SCROLL_DOWN row (now row 14 after insertion) to TINY_SCROLL_DOWN
$ws.Cells.Item(14, 2).Value = "TINY_SCROLL_DOWN"

# Delete the "ViewFullCart" row (now row 19 after insertion)
$ws.Rows.Item(19).Delete()

$wb.Save()
